# Update the search item in cell A2 from "coffee mug" to "Aleheida"
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = "Aleheida"

# Reflect the active selection change (was B8, now A2) recorded in the saved workbook
$ws.Range("A2").Select()
